$wb = $excel.ActiveWorkbook

# Sheet 1: "OFF"
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 273
$wsOff.Range("C3").Value = 206
$wsOff.Range("D3").Value = 65

# Sheet 2: "DEF"
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 239
$wsDef.Range("C3").Value = 145
$wsDef.Range("D3").Value = 63
$wsDef.Range("E3").Value = 26
$wsDef.Range("F3").Value = 9
$wsDef.Range("G3").Value = 6
